# Auto-update data + news
# Updates macro credit metrics sheet: refresh value/as_of/avg10/delta columns
# for the "Card 30+ Delinquency" and "Net Charge-off Rate" rows, flip the
# Net Charge-off Rate status back to "healthy", and refresh the
# "Mortgage 30+ Delinquency" row's as_of/avg10/delta columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: DRCCLACBS_pct (Card 30+ Delinquency)
$ws.Range("E2").Value = 2.94
$ws.Range("F2").Value = "'Oct 2025"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 2.465121951219512
$ws.Range("H2").Value = -0.1400000000000001
$ws.Range("I2").Value = -0.04545454545454549

# Row 3: CORCCACBS_pct (Net Charge-off Rate)
$ws.Range("D3").Value = "healthy"
$ws.Range("E3").Value = 4.11
$ws.Range("F3").Value = "'Oct 2025"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 3.390243902439025
$ws.Range("H3").Value = -0.4699999999999998
$ws.Range("I3").Value = -0.1026200873362445

# Row 5: DRSFRMACBS_pct (Mortgage 30+ Delinquency)
$ws.Range("F5").Value = "'Oct 2025"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = 2.664634146341463
$ws.Range("H5").Value = 0.01000000000000001
$ws.Range("I5").Value = 0.005649717514124299
